$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row
$ws.Range("A1").Value = "nome"
$ws.Range("B1").Value = "peso"
$ws.Range("C1").Value = "zone"
$ws.Range("D1").Value = "partenza"

# Data row 2
$ws.Range("A2").Value = "40279 (interno)"
$ws.Range("B2").Value = "CAMPO VUOTO"
$ws.Range("C2").Value = "CAMPO VUOTO"
$ws.Range("D2").Value = 0

# Data row 3
$ws.Range("A3").Value = "40176 (interno)"
$ws.Range("B3").Value = "CAMPO VUOTO"
$ws.Range("C3").Value = "CAMPO VUOTO"
$ws.Range("D3").Value = 0

# Highlight the "CAMPO VUOTO" cells with a yellow fill
$ws.Range("B2:C3").Interior.Color = 65535

# Set column widths (30 chars) for columns A through D
$ws.Range("A1:D1").EntireColumn.ColumnWidth = 29.166666666666668
